$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data appended after the existing table (rows 2-5), representing
# freshly executed "small" problem results from the maps/json folder.

$rows = @(
    @{ A = "avenida_de_espanÌa_250_0.json"; B = "DFS"; C = 9;  D = 14; E = 4;  F = 6.46256666666666;    G = 20; H = "small" },
    @{ A = "avenida_de_espanÌa_250_0.json"; B = "BFS"; C = 10; D = 14; E = 3;  F = 4.23893333333333;    G = 42; H = "small" },
    @{ A = "avenida_de_espanÌa_250_1.json "; B = "DFS"; C = 14; D = 18; E = 9;  F = 15.4287333333333;   G = 25; H = "small" },
    @{ A = "avenida_de_espanÌa_250_1.json"; B = "BFS"; C = 24; D = 39; E = 9;  F = 15.4287333333333;   G = 70; H = "small" },
    @{ A = "paseo_simoÌn_abril_250_0.json"; B = "DFS" },
    @{ A = "paseo_simoÌn_abril_250_0.json"; B = "BFS" },
    @{ A = "paseo_simoÌn_abril_250_1.json"; B = "DFS"; C = 25; D = 36; E = 12; F = 27.5196166666666;   G = 43; H = "small" },
    @{ A = "paseo_simoÌn_abril_250_1.json"; B = "BFS"; C = 24; D = 35; E = 8;  F = 20.1584333333333;   G = 84; H = "small" }
)

$startRow = 6
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    if ($row.ContainsKey("C")) {
        $ws.Cells.Item($r, 3).Value = $row.C
        $ws.Cells.Item($r, 4).Value = $row.D
        $ws.Cells.Item($r, 5).Value = $row.E
        $ws.Cells.Item($r, 6).Value = $row.F
        $ws.Cells.Item($r, 7).Value = $row.G
        $ws.Cells.Item($r, 8).Value = $row.H
    }
}

$ws.Range("I13").Select()
